$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Reorder rows 7-9.
#    Before: row7=A 34145-2019, row8=A 37380-2020, row9=A 66048-2020 (stale)
#    After : row7=A 66048-2020 (updated figures), row8=A 34145-2019, row9=A 37380-2020
#    Achieved by deleting the stale row 9 and inserting a fresh row at 7,
#    which leaves rows 7/8 (34145/37380) shifted down by one into 8/9,
#    and a blank row 7 ready to receive the updated data.
# ------------------------------------------------------------------
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Insert()

# ------------------------------------------------------------------
# 2) Fill in the new row 7 with the refreshed "A 66048-2020" record.
# ------------------------------------------------------------------
$ws.Range("A7").Value = "A 66048-2020"
$ws.Range("B7").Value = 44175
$ws.Range("B7").NumberFormat = $ws.Range("B8").NumberFormat
$ws.Range("C7").Value = 45202
$ws.Range("C7").NumberFormat = $ws.Range("C8").NumberFormat
$ws.Range("D7").Value = "SÖDERMANLANDS LÄN"
$ws.Range("E7").Value = "STRÄNGNÄS"
$ws.Range("F7").Value = "Övriga Aktiebolag"
$ws.Range("G7").Value = 4.6
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 4
$ws.Range("R7").Value = "Grönpyrola`r`nKantarellvaxskivling`r`nPurpurknipprot`r`nTibast"
$ws.Range("R7").WrapText = $true

$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 66048-2020.xlsx", "A 66048-2020")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 66048-2020.png", "A 66048-2020")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 66048-2020.docx", "A 66048-2020")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 66048-2020.docx", "A 66048-2020")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 66048-2020.docx", "A 66048-2020")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 66048-2020.docx", "A 66048-2020")'

$ws.Rows.Item(7).RowHeight = 15

# ------------------------------------------------------------------
# 3) Update "Förändrad" (column C) for every existing data row (2-432)
#    from 45192 to 45202.
# ------------------------------------------------------------------
for ($r = 2; $r -le 432; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 432 did not previously carry an explicit row height; it now does.
$ws.Rows.Item(432).RowHeight = 15

# ------------------------------------------------------------------
# 4) Append the new record "A 45742-2023" as row 433.
# ------------------------------------------------------------------
$ws.Range("A433").Value = "A 45742-2023"
$ws.Range("B433").Value = 45195
$ws.Range("B433").NumberFormat = $ws.Range("B432").NumberFormat
$ws.Range("C433").Value = 45202
$ws.Range("C433").NumberFormat = $ws.Range("C432").NumberFormat
$ws.Range("D433").Value = "SÖDERMANLANDS LÄN"
$ws.Range("E433").Value = "STRÄNGNÄS"
$ws.Range("F433").Value = "Sveaskog"
$ws.Range("G433").Value = 1.8
$ws.Range("H433").Value = 0
$ws.Range("I433").Value = 0
$ws.Range("J433").Value = 0
$ws.Range("K433").Value = 0
$ws.Range("L433").Value = 0
$ws.Range("M433").Value = 0
$ws.Range("N433").Value = 0
$ws.Range("O433").Value = 0
$ws.Range("P433").Value = 0
$ws.Range("Q433").Value = 0
$ws.Range("R433").Value = ""
$ws.Range("R433").WrapText = $true
